# Updated cryptos list on Mon Dec  2 13:03:32 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.008.67"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "3.595.39"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +23.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "222.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "633.15"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.411"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "3.591.66"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.206"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000287"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("D16").Value = "4.265.20"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "94.807.17"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.78%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.49%  "
$ws.Range("B21").Value = "WrappedEther"
$ws.Range("C21").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D21").Value = "3.590.74"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.515"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "501.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.27%  "
$ws.Range("E25").Value = "  +21.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "120.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +18.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000200"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("D29").Value = "3.786.70"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  -5.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.47%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "584.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.157"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.475"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0471"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.916"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("B50").Value = "MantraDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "220.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.31%  "
